$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card7")

# --- Header row (row 1) ---
# Fix existing "Correction " header text: drop the trailing space.
$ws.Range("N1").Value = "Correction"

# Copy N1's header formatting (bold font, border, center/top alignment) onto
# the brand new O1 header cell before filling in its text.
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("O1").Value = "Serviced by "

# --- Data rows (2-12) ---
# Column N used to hold empty placeholder cells; they now get the same "nan"
# placeholder text used throughout the rest of the row.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

# Column O is brand new. Every data row needs an (empty) placeholder cell so
# the column physically exists across the whole table, matching the other
# still-blank columns. Touch each cell (border on/off) so Excel materializes
# it in the sheet instead of treating it as untouched.
for ($r = 2; $r -le 12; $r++) {
    $c = $ws.Cells.Item($r, 15)
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(7).LineStyle = -4142
}
